$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.548.65"
$ws.Range("E2").Value = "'  -1.20%  "
$ws.Range("D3").Value = "'3.908.10"
$ws.Range("E3").Value = "'  +3.11%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'599.90"
$ws.Range("E5").Value = "'  -0.49%  "
$ws.Range("D6").Value = "'164.47"
$ws.Range("E6").Value = "'  -1.24%  "
$ws.Range("D7").Value = "'3.906.19"
$ws.Range("E7").Value = "'  +3.13%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "'  -2.47%  "
$ws.Range("E10").Value = "'  -4.29%  "
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "'  -0.03%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "'  -0.83%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "'  -1.15%  "
$ws.Range("D14").Value = "'36.81"
$ws.Range("E14").Value = "'  -2.75%  "
$ws.Range("D15").Value = "'4.561.36"
$ws.Range("E15").Value = "'  +3.15%  "
$ws.Range("D16").Value = "'3.929.94"
$ws.Range("E16").Value = "'  +3.87%  "
$ws.Range("D17").Value = "'68.718.49"
$ws.Range("E17").Value = "'  -1.03%  "
$ws.Range("E18").Value = "'  -0.82%  "
$ws.Range("E19").Value = "'  -1.68%  "
$ws.Range("D20").Value = "'17.00"
$ws.Range("E20").Value = "'  -3.84%  "
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "'  -1.87%  "
$ws.Range("D22").Value = "'483.58"
$ws.Range("E22").Value = "'  -2.29%  "
$ws.Range("B23").Value = "'Polygon"
$ws.Range("C23").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "'  -1.48%  "
$ws.Range("B24").Value = "'PEPE"
$ws.Range("C24").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "'0.0000168"
$ws.Range("E24").Value = "'  +10.68%  "
$ws.Range("D25").Value = "'84.34"
$ws.Range("E25").Value = "'  -0.82%  "
$ws.Range("E26").Value = "'  -2.19%  "
$ws.Range("D27").Value = "'11.96"
$ws.Range("E27").Value = "'  -2.97%  "
$ws.Range("E28").Value = "'  -0.92%  "
$ws.Range("E29").Value = "'  -0.15%  "
$ws.Range("E30").Value = "'  -1.80%  "
$ws.Range("D31").Value = "'4.058.02"
$ws.Range("E31").Value = "'  +3.40%  "
$ws.Range("D32").Value = "'7.81"
$ws.Range("E32").Value = "'  -4.22%  "
$ws.Range("E33").Value = "'  -3.02%  "
$ws.Range("D34").Value = "'31.96"
$ws.Range("E34").Value = "'  -0.36%  "
$ws.Range("D35").Value = "'3.852.11"
$ws.Range("E35").Value = "'  +3.14%  "
$ws.Range("E36").Value = "'  -1.79%  "
$ws.Range("D37").Value = "'1.03"
$ws.Range("E37").Value = "'  +1.65%  "
$ws.Range("E38").Value = "'  -0.25%  "
$ws.Range("D39").Value = "'5.87"
$ws.Range("E39").Value = "'  -1.79%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "'  +0.27%  "
$ws.Range("D41").Value = "'3.07"
$ws.Range("E41").Value = "'  -0.29%  "
$ws.Range("E42").Value = "'  -3.14%  "
$ws.Range("D43").Value = "'433.26"
$ws.Range("E43").Value = "'  +1.40%  "
$ws.Range("D44").Value = "'48.48"
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("D45").Value = "'1.98"
$ws.Range("E45").Value = "'  -1.09%  "
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("E47").Value = "'  -0.79%  "
$ws.Range("D48").Value = "'26.45"
$ws.Range("E48").Value = "'  +9.99%  "
$ws.Range("D49").Value = "'141.93"
$ws.Range("E49").Value = "'  +0.40%  "
$ws.Range("D50").Value = "'2.814.27"
$ws.Range("E50").Value = "'  -0.10%  "
$ws.Range("B51").Value = "'Arweave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'39.24"
$ws.Range("E51").Value = "'  -2.52%  "
